# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp: 08:22 -> 08:52
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 08:52"

# Singapur (row 31): new-cases (F) 26 -> 24
$ws.Range("F31").Value = 24

# Ucrania (row 41): refreshed daily figures
$ws.Range("B41").Value = 8125
$ws.Range("C41").Value = 478
$ws.Range("D41").Value = 782
$ws.Range("E41").Value = 7142
$ws.Range("F41").Value = 104
$ws.Range("G41").Value = 8
$ws.Range("H41").Value = 201

# Chequia (row 44): refreshed daily figures (B/C unchanged)
$ws.Range("D44").Value = 2389
$ws.Range("E44").Value = 4669
$ws.Range("F44").Value = 71
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 215

# Malta / Georgia swap places (Georgia overtakes Malta in total cases)
# Row 107 now holds Georgia's refreshed figures
$ws.Range("A107").Value = "Georgia"
$ws.Range("B107").Value = 456
$ws.Range("C107").Value = 12
$ws.Range("D107").Value = 132
$ws.Range("E107").Value = 319
$ws.Range("F107").Value = 6
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 5

# Row 108 now holds Malta's (previous) figures
$ws.Range("A108").Value = "Malta"
$ws.Range("B108").Value = 447
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 223
$ws.Range("E108").Value = 221
$ws.Range("F108").Value = 2
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 3

# Taiwan (row 111): refreshed daily figures
$ws.Range("B111").Value = 429
$ws.Range("C111").Value = 1
$ws.Range("D111").Value = 275
$ws.Range("E111").Value = 148
